$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.584.34'
$ws.Range('E2').Value = '  +5.99%  '
$ws.Range('D3').Value = '2.338.34'
$ws.Range('E3').Value = '  +2.93%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '0.997'
$c.ClearFormats()
$ws.Range('E4').Value = '  -0.52%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '521.96'
$c.ClearFormats()
$ws.Range('E5').Value = '  +4.46%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '135.39'
$c.ClearFormats()
$ws.Range('E6').Value = '  +4.83%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.993'
$c.ClearFormats()
$ws.Range('E7').Value = '  -0.48%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.539'
$c.ClearFormats()
$ws.Range('E8').Value = '  +2.44%  '
$ws.Range('D9').Value = '2.370.85'
$ws.Range('E9').Value = '  +4.03%  '
$ws.Range('E10').Value = '  +9.20%  '
$ws.Range('E11').Value = '  +0.80%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '5.23'
$c.ClearFormats()
$ws.Range('E12').Value = '  +6.74%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '0.346'
$c.ClearFormats()
$ws.Range('E13').Value = '  +2.63%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '24.03'
$c.ClearFormats()
$ws.Range('E14').Value = '  +3.63%  '
$ws.Range('D15').Value = '2.758.31'
$ws.Range('E15').Value = '  +3.17%  '
$ws.Range('D16').Value = '56.972.95'
$ws.Range('E16').Value = '  +4.87%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.0000136'
$c.ClearFormats()
$ws.Range('E17').Value = '  +5.05%  '
$ws.Range('D18').Value = '2.360.61'
$ws.Range('E18').Value = '  +3.45%  '
$ws.Range('E19').Value = '  +3.16%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '4.31'
$c.ClearFormats()
$ws.Range('E20').Value = '  +3.67%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '323.35'
$c.ClearFormats()
$ws.Range('E21').Value = '  +6.26%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '6.71'
$c.ClearFormats()
$ws.Range('E22').Value = '  +6.34%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.ClearFormats()
$ws.Range('E23').Value = '  +0.01%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '61.33'
$c.ClearFormats()
$ws.Range('E24').Value = '  +1.78%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '0.161'
$c.ClearFormats()
$ws.Range('E25').Value = '  +7.92%  '
$ws.Range('E26').Value = '  -0.85%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '7.84'
$c.ClearFormats()
$ws.Range('E27').Value = '  +6.32%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '171.99'
$c.ClearFormats()
$ws.Range('E28').Value = '  -1.94%  '
$ws.Range('D29').Value = '0.0₃0747'
$ws.Range('E29').Value = '  +5.95%  '
$ws.Range('E30').Value = '  +11.63%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '6.35'
$c.ClearFormats()
$ws.Range('E31').Value = '  +5.52%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '18.49'
$c.ClearFormats()
$ws.Range('E33').Value = '  +3.64%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.998'
$c.ClearFormats()
$ws.Range('E34').Value = '  -0.04%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '0.960'
$c.ClearFormats()
$ws.Range('E35').Value = '  +0.94%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.995'
$c.ClearFormats()
$ws.Range('E36').Value = '  -0.20%  '
$ws.Range('E37').Value = '  +5.74%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '4.07'
$c.ClearFormats()
$ws.Range('E38').Value = '  +8.93%  '
$ws.Range('E39').Value = '  +9.00%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '37.63'
$c.ClearFormats()
$ws.Range('E40').Value = '  +4.14%  '
$ws.Range('E41').Value = '  +2.48%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '140.57'
$c.ClearFormats()
$ws.Range('E42').Value = '  +12.58%  '
$ws.Range('E43').Value = '  +6.85%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '5.28'
$c.ClearFormats()
$ws.Range('E44').Value = '  +9.34%  '
$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '281.99'
$c.ClearFormats()
$ws.Range('E45').Value = '  +14.33%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '0.0513'
$c.ClearFormats()
$ws.Range('E46').Value = '  +4.17%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.0933'
$c.ClearFormats()
$ws.Range('E47').Value = '  +3.72%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '0.567'
$c.ClearFormats()
$ws.Range('E48').Value = '  +3.75%  '
$ws.Range('B49').Value = 'Polygon'
$ws.Range('C49').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.384'
$c.ClearFormats()
$ws.Range('E49').Value = '  +2.62%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0217'
$c.ClearFormats()
$ws.Range('E50').Value = '  +5.70%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '17.09'
$c.ClearFormats()
$ws.Range('E51').Value = '  +4.89%  '
